# Rotate the "Recorded By" (column G) author list so the first entry moves
# to the end, e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System".
# Cells with only a single entry (no comma) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $value = $cell.Value2

    if ($value -ne $null) {
        $parts = $value -split ", "

        if ($parts.Count -gt 1) {
            $rotated = $parts[1..($parts.Count - 1)] + $parts[0]
            $cell.Value = $rotated -join ", "
        }
    }
}
